$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = [double]"3"
$ws.Range("F2").Value = [double]"1"
$ws.Range("G2").Value = [double]"0.09665433333333333"
$ws.Range("H2").Value = [double]"0.289963"
$ws.Range("I2").Value = [double]"0.0006230336790718351"
$ws.Range("J2").Value = [double]"0.0006230336790718351"
$ws.Range("M2").Value = [double]"2.027115333333333"
$ws.Range("N2").Value = [double]"6.081346"
$ws.Range("O2").Value = [double]"0.006596284565418616"
$ws.Range("P2").Value = [double]"0.006596284565418615"
$ws.Range("Q2").Value = [double]"0.1959294811331111"
$ws.Range("R2").Value = [double]"1.763365330198"
$ws.Range("S2").Value = [double]"4.109707440997521E-06"
$ws.Range("T2").Value = [double]"4.10970744099752E-06"

# Row 3
$ws.Range("E3").Value = [double]"3"
$ws.Range("F3").Value = [double]"1"
$ws.Range("G3").Value = [double]"0.09665433333333333"
$ws.Range("H3").Value = [double]"0.289963"
$ws.Range("I3").Value = [double]"0.0006230336790718351"
$ws.Range("J3").Value = [double]"0.0006230336790718351"
$ws.Range("M3").Value = [double]"256.4443053333333"
$ws.Range("N3").Value = [double]"769.332916"
$ws.Range("O3").Value = [double]"0.8344762556643375"
$ws.Range("P3").Value = [double]"0.8344762556643374"
$ws.Range("Q3").Value = [double]"24.78645336912311"
$ws.Range("R3").Value = [double]"223.078080322108"
$ws.Range("S3").Value = [double]"0.0005199068116646414"
$ws.Range("T3").Value = [double]"0.0005199068116646413"

# Row 4
$ws.Range("E4").Value = [double]"3"
$ws.Range("F4").Value = [double]"1"
$ws.Range("G4").Value = [double]"0.09665433333333333"
$ws.Range("H4").Value = [double]"0.289963"
$ws.Range("I4").Value = [double]"0.0006230336790718351"
$ws.Range("J4").Value = [double]"0.0006230336790718351"
$ws.Range("M4").Value = [double]"48.84026566666667"
$ws.Range("N4").Value = [double]"146.520797"
$ws.Range("O4").Value = [double]"0.158927459770244"
$ws.Range("P4").Value = [double]"0.158927459770244"
$ws.Range("Q4").Value = [double]"4.720623317834556"
$ws.Range("R4").Value = [double]"42.485609860511"
$ws.Range("S4").Value = [double]"9.901715996619619E-05"
$ws.Range("T4").Value = [double]"9.901715996619616E-05"

# Row 5
$ws.Range("E5").Value = [double]"3"
$ws.Range("F5").Value = [double]"1"
$ws.Range("G5").Value = [double]"154.8642143333334"
$ws.Range("H5").Value = [double]"464.5926430000001"
$ws.Range("I5").Value = [double]"0.9982544794956518"
$ws.Range("J5").Value = [double]"0.9982544794956519"
$ws.Range("M5").Value = [double]"2.027115333333333"
$ws.Range("N5").Value = [double]"6.081346"
$ws.Range("O5").Value = [double]"0.006596284565418616"
$ws.Range("P5").Value = [double]"0.006596284565418615"
$ws.Range("Q5").Value = [double]"313.9276234597198"
$ws.Range("R5").Value = [double]"2825.348611137478"
$ws.Range("S5").Value = [double]"0.006584770615457162"
$ws.Range("T5").Value = [double]"0.006584770615457162"

# Row 6
$ws.Range("E6").Value = [double]"3"
$ws.Range("F6").Value = [double]"1"
$ws.Range("G6").Value = [double]"154.8642143333334"
$ws.Range("H6").Value = [double]"464.5926430000001"
$ws.Range("I6").Value = [double]"0.9982544794956518"
$ws.Range("J6").Value = [double]"0.9982544794956519"
$ws.Range("M6").Value = [double]"256.4443053333333"
$ws.Range("N6").Value = [double]"769.332916"
$ws.Range("O6").Value = [double]"0.8344762556643375"
$ws.Range("P6").Value = [double]"0.8344762556643374"
$ws.Range("Q6").Value = [double]"39714.04586570412"
$ws.Range("R6").Value = [double]"357426.412791337"
$ws.Range("S6").Value = [double]"0.8330196602496837"
$ws.Range("T6").Value = [double]"0.8330196602496837"

# Row 7
$ws.Range("E7").Value = [double]"3"
$ws.Range("F7").Value = [double]"1"
$ws.Range("G7").Value = [double]"154.8642143333334"
$ws.Range("H7").Value = [double]"464.5926430000001"
$ws.Range("I7").Value = [double]"0.9982544794956518"
$ws.Range("J7").Value = [double]"0.9982544794956519"
$ws.Range("M7").Value = [double]"48.84026566666667"
$ws.Range("N7").Value = [double]"146.520797"
$ws.Range("O7").Value = [double]"0.158927459770244"
$ws.Range("P7").Value = [double]"0.158927459770244"
$ws.Range("Q7").Value = [double]"7563.60937029961"
$ws.Range("R7").Value = [double]"68072.48433269649"
$ws.Range("S7").Value = [double]"0.1586500486305111"
$ws.Range("T7").Value = [double]"0.158650048630511"

# Row 8
$ws.Range("E8").Value = [double]"3"
$ws.Range("F8").Value = [double]"1"
$ws.Range("G8").Value = [double]"0.174137"
$ws.Range("H8").Value = [double]"0.522411"
$ws.Range("I8").Value = [double]"0.001122486825276316"
$ws.Range("J8").Value = [double]"0.001122486825276316"
$ws.Range("M8").Value = [double]"2.027115333333333"
$ws.Range("N8").Value = [double]"6.081346"
$ws.Range("O8").Value = [double]"0.006596284565418616"
$ws.Range("P8").Value = [double]"0.006596284565418615"
$ws.Range("Q8").Value = [double]"0.3529957828006666"
$ws.Range("R8").Value = [double]"3.176962045206"
$ws.Range("S8").Value = [double]"7.404242520455905E-06"
$ws.Range("T8").Value = [double]"7.404242520455905E-06"

# Row 9
$ws.Range("E9").Value = [double]"3"
$ws.Range("F9").Value = [double]"1"
$ws.Range("G9").Value = [double]"0.174137"
$ws.Range("H9").Value = [double]"0.522411"
$ws.Range("I9").Value = [double]"0.001122486825276316"
$ws.Range("J9").Value = [double]"0.001122486825276316"
$ws.Range("M9").Value = [double]"256.4443053333333"
$ws.Range("N9").Value = [double]"769.332916"
$ws.Range("O9").Value = [double]"0.8344762556643375"
$ws.Range("P9").Value = [double]"0.8344762556643374"
$ws.Range("Q9").Value = [double]"44.65644199783066"
$ws.Range("R9").Value = [double]"401.9079779804759"
$ws.Range("S9").Value = [double]"0.0009366886029891294"
$ws.Range("T9").Value = [double]"0.0009366886029891295"

# Row 10
$ws.Range("E10").Value = [double]"3"
$ws.Range("F10").Value = [double]"1"
$ws.Range("G10").Value = [double]"0.174137"
$ws.Range("H10").Value = [double]"0.522411"
$ws.Range("I10").Value = [double]"0.001122486825276316"
$ws.Range("J10").Value = [double]"0.001122486825276316"
$ws.Range("M10").Value = [double]"48.84026566666667"
$ws.Range("N10").Value = [double]"146.520797"
$ws.Range("O10").Value = [double]"0.158927459770244"
$ws.Range("P10").Value = [double]"0.158927459770244"
$ws.Range("Q10").Value = [double]"8.504897342396333"
$ws.Range("R10").Value = [double]"76.544076081567"
$ws.Range("S10").Value = [double]"0.0001783939797667306"
$ws.Range("T10").Value = [double]"0.0001783939797667306"
